$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# The activity "Création des uses cases et scénarios" originally had "4 heures "
# logged in C5; update it to "5 heures " (trailing space preserved).
$ws.Range("C5").Value = "5 heures "
